# Apply cryptos list update (price/volume refresh) per commit
# "Updated cryptos list on Thu Jan 11 13:28:26 UTC 2024 with GitHub Actions"
#
# Each cell write uses a leading apostrophe to force text entry (so
# numeric-looking values like "315.56" or "0.620" are stored verbatim
# as text instead of being coerced into Excel numbers), then resets the
# cell Style back to "Normal" to drop the transient quote-prefix marker
# so the saved style matches the original (unstyled) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) / Volume(1h) (E) refresh for existing rows ---
$ws.Range("D2").Value = "'" + '47.317.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +4.98%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '2.653.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +11.48%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + '  +0.14%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '315.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  +7.68%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '106.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +12.63%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'" + '0.620'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + '  +11.41%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'" + '  +0.03%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'" + '0.608'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + '  +22.28%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + '40.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  +18.05%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + '56.05'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  +5.00%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'" + '0.0863'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + '  +11.65%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '8.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +23.62%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + '3.059.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +11.53%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'" + '  +3.43%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '2.671.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  +12.43%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '0.955'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  +16.40%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'" + '15.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + '  +12.63%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + '48.114.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  +6.95%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '0.0000105'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  +12.86%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + '13.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +8.97%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + '6.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  +14.29%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + '73.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  +10.51%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'" + '273.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + '  +14.77%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + '3.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + '  +15.16%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '31.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +48.41%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + '2.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +19.03%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'" + '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + '  +0.26%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'" + '4.08'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  +1.48%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '10.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  +14.21%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'" + '40.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +7.77%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'" + '  +5.02%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'" + '6.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + '  +17.85%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + '3.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +0.04%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'" + '2.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +18.74%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D38").Value = "'" + '153.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  +4.55%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'" + '  +10.76%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + '0.126'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  +10.26%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + '24.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +62.83%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + '17.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  +16.01%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'" + '4.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  +16.88%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'" + '3.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +20.39%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + '0.0338'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  +15.52%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + '2.234.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +15.98%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '97.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  +8.93%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D50").Value = "'" + '116.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +16.66%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'" + '  +9.60%  '
$ws.Range("E51").Style = "Normal"

# --- Row 36/37 swap: Hedera now ranks above WEMIXToken ---
$ws.Range("B36").Value = "'" + 'Hedera'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'" + '0.0864'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +14.25%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'" + 'WEMIXToken'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'" + 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'" + '2.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  +7.29%  '
$ws.Range("E37").Style = "Normal"

# --- Row 48/49 swap: FraxShare now ranks above FirstDigitalUSD ---
$ws.Range("B48").Value = "'" + 'FraxShare'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'" + 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'" + '10.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  +22.84%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'" + 'FirstDigitalUSD'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'" + 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'" + '0.999'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  +0.03%  '
$ws.Range("E49").Style = "Normal"
